$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dCell = $ws.Range("D2")
$dCell.NumberFormat = "@"
$dCell.Value = "43.232.13"
$dCell.Style = "Normal"
$ws.Range("E2").Value = "  -4.63%  "
$dCell = $ws.Range("D3")
$dCell.NumberFormat = "@"
$dCell.Value = "2.239.15"
$dCell.Style = "Normal"
$ws.Range("E3").Value = "  -5.55%  "
$ws.Range("E4").Value = "  +0.03%  "
$dCell = $ws.Range("D5")
$dCell.NumberFormat = "@"
$dCell.Value = "319.11"
$dCell.Style = "Normal"
$ws.Range("E5").Value = "  +1.38%  "
$dCell = $ws.Range("D6")
$dCell.NumberFormat = "@"
$dCell.Value = "100.95"
$dCell.Style = "Normal"
$ws.Range("E6").Value = "  -6.57%  "
$dCell = $ws.Range("D7")
$dCell.NumberFormat = "@"
$dCell.Value = "0.587"
$dCell.Style = "Normal"
$ws.Range("E7").Value = "  -8.05%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("E9").Value = "  -7.92%  "
$dCell = $ws.Range("D10")
$dCell.NumberFormat = "@"
$dCell.Value = "37.10"
$dCell.Style = "Normal"
$ws.Range("E10").Value = "  -9.08%  "
$dCell = $ws.Range("D11")
$dCell.NumberFormat = "@"
$dCell.Value = "54.50"
$dCell.Style = "Normal"
$ws.Range("E11").Value = "  -2.51%  "
$dCell = $ws.Range("D12")
$dCell.NumberFormat = "@"
$dCell.Value = "0.0826"
$dCell.Style = "Normal"
$ws.Range("E12").Value = "  -9.89%  "
$dCell = $ws.Range("D13")
$dCell.NumberFormat = "@"
$dCell.Value = "7.71"
$dCell.Style = "Normal"
$ws.Range("E13").Value = "  -9.16%  "
$ws.Range("E14").Value = "  -1.10%  "
$dCell = $ws.Range("D15")
$dCell.NumberFormat = "@"
$dCell.Value = "2.578.56"
$dCell.Style = "Normal"
$ws.Range("E15").Value = "  -5.57%  "
$dCell = $ws.Range("D16")
$dCell.NumberFormat = "@"
$dCell.Value = "0.866"
$dCell.Style = "Normal"
$ws.Range("E16").Value = "  -11.64%  "
$dCell = $ws.Range("D17")
$dCell.NumberFormat = "@"
$dCell.Value = "14.41"
$dCell.Style = "Normal"
$ws.Range("E17").Value = "  -6.08%  "
$dCell = $ws.Range("D18")
$dCell.NumberFormat = "@"
$dCell.Value = "2.238.81"
$dCell.Style = "Normal"
$ws.Range("E18").Value = "  -5.50%  "
$dCell = $ws.Range("D19")
$dCell.NumberFormat = "@"
$dCell.Value = "43.168.51"
$dCell.Style = "Normal"
$ws.Range("E19").Value = "  -4.72%  "
$dCell = $ws.Range("D20")
$dCell.NumberFormat = "@"
$dCell.Value = "14.30"
$dCell.Style = "Normal"
$ws.Range("E20").Value = "  -8.41%  "
$dCell = $ws.Range("D21")
$dCell.NumberFormat = "@"
$dCell.Value = "0.0₃0967"
$dCell.Style = "Normal"
$ws.Range("E21").Value = "  -8.74%  "
$ws.Range("E22").Value = "  -9.98%  "
$dCell = $ws.Range("D23")
$dCell.NumberFormat = "@"
$dCell.Value = "65.54"
$dCell.Style = "Normal"
$ws.Range("E23").Value = "  -10.53%  "
$ws.Range("E24").Value = "  -11.66%  "
$dCell = $ws.Range("D25")
$dCell.NumberFormat = "@"
$dCell.Value = "238.31"
$dCell.Style = "Normal"
$ws.Range("E25").Value = "  -8.60%  "
$ws.Range("E26").Value = "  -7.98%  "
$dCell = $ws.Range("D27")
$dCell.NumberFormat = "@"
$dCell.Value = "1.00"
$dCell.Style = "Normal"
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("E28").Value = "  +2.25%  "
$dCell = $ws.Range("D29")
$dCell.NumberFormat = "@"
$dCell.Value = "10.08"
$dCell.Style = "Normal"
$ws.Range("E29").Value = "  -9.67%  "
$ws.Range("E30").Value = "  -2.41%  "
$ws.Range("E31").Value = "  -14.48%  "
$dCell = $ws.Range("D32")
$dCell.NumberFormat = "@"
$dCell.Value = "35.42"
$dCell.Style = "Normal"
$ws.Range("E32").Value = "  -4.30%  "
$dCell = $ws.Range("D33")
$dCell.NumberFormat = "@"
$dCell.Value = "20.51"
$dCell.Style = "Normal"
$ws.Range("E33").Value = "  -8.09%  "
$dCell = $ws.Range("D34")
$dCell.NumberFormat = "@"
$dCell.Value = "0.0878"
$dCell.Style = "Normal"
$ws.Range("E34").Value = "  -9.40%  "
$dCell = $ws.Range("D35")
$dCell.NumberFormat = "@"
$dCell.Value = "153.60"
$dCell.Style = "Normal"
$ws.Range("E35").Value = "  -7.90%  "
$ws.Range("E36").Value = "  -3.62%  "
$dCell = $ws.Range("D37")
$dCell.NumberFormat = "@"
$dCell.Value = "3.20"
$dCell.Style = "Normal"
$ws.Range("E37").Value = "  +9.64%  "
$ws.Range("E38").Value = "  +5.28%  "
$ws.Range("E39").Value = "  -7.33%  "
$ws.Range("E40").Value = "  -4.61%  "
$ws.Range("E41").Value = "  -10.72%  "
$dCell = $ws.Range("D42")
$dCell.NumberFormat = "@"
$dCell.Value = "3.70"
$dCell.Style = "Normal"
$ws.Range("E42").Value = "  -6.84%  "
$dCell = $ws.Range("D43")
$dCell.NumberFormat = "@"
$dCell.Value = "0.0324"
$dCell.Style = "Normal"
$ws.Range("E43").Value = "  -8.48%  "
$dCell = $ws.Range("D44")
$dCell.NumberFormat = "@"
$dCell.Value = "12.86"
$dCell.Style = "Normal"
$ws.Range("E44").Value = "  -1.62%  "
$ws.Range("E45").Value = "  -0.02%  "
$dCell = $ws.Range("D46")
$dCell.NumberFormat = "@"
$dCell.Value = "1.794.43"
$dCell.Style = "Normal"
$ws.Range("E46").Value = "  -1.25%  "
$dCell = $ws.Range("D47")
$dCell.NumberFormat = "@"
$dCell.Value = "87.19"
$dCell.Style = "Normal"
$ws.Range("E47").Value = "  -11.36%  "
$ws.Range("E48").Value = "  -9.55%  "
$dCell = $ws.Range("D49")
$dCell.NumberFormat = "@"
$dCell.Value = "76.67"
$dCell.Style = "Normal"
$ws.Range("E49").Value = "  -7.78%  "
$dCell = $ws.Range("D50")
$dCell.NumberFormat = "@"
$dCell.Value = "5.33"
$dCell.Style = "Normal"
$ws.Range("E50").Value = "  -9.96%  "
$dCell = $ws.Range("D51")
$dCell.NumberFormat = "@"
$dCell.Value = "59.59"
$dCell.Style = "Normal"
$ws.Range("E51").Value = "  -15.06%  "
